# Edit sheet Card24 by admin
# The 3rd data row (worksheet row 3, "0 / 150 / 99 / done / done / ... / تم معايره")
# was removed from the Card24 lookup table; all subsequent rows shift up by
# one and the last row (former row 13) disappears, shrinking the used range
# from A1:N13 to A1:N12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

# Deleting the entire row shifts rows 4-13 up to 3-12 and updates the
# worksheet dimension automatically.
$ws.Rows(3).Delete()
